# Fix Ubuntu AMI IDs on the pricing/invoice sheet (Sheet1), row 4 ("Ubuntu 14").
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Columns E, G, H, I, J hold the existing per-region AMI ids for Ubuntu 14;
# replace the stale ids with the corrected ones, and add a new region (K4).
$ws.Range("E4").Value = "ami-9dde7f8b"
$ws.Range("G4").Value = "ami-9d772efd"
$ws.Range("H4").Value = "ami-0e2aa66e"
$ws.Range("I4").Value = "ami-115d7777"
$ws.Range("J4").Value = "ami-6039ed0f"
$ws.Range("K4").Value = "ami-c29184a6"

# Match the author's resulting selection state in the saved file.
$ws.Range("L22").Select()
